$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 652.84705
$ws.Range("J17").Value = 549.1594
$ws.Range("L17").Value = 1647.4782
$ws.Range("N17").Value = -1983.4782
$ws.Range("H38").Value = 3860.6155
$ws.Range("I38").Value = 114.5
$ws.Range("J38").Value = 7071.5713
$ws.Range("K38").Value = 343.5
$ws.Range("L38").Value = 21214.7139
$ws.Range("M38").Value = 28.5
$ws.Range("N38").Value = -21958.7139
$ws.Range("H86").Value = 1550.4
$ws.Range("I86").Value = 1481
$ws.Range("J86").Value = 1679.2858
$ws.Range("K86").Value = 1481
$ws.Range("L86").Value = 1679.2858
$ws.Range("M86").Value = -358
$ws.Range("N86").Value = -3925.2858
$ws.Range("H89").Value = 1550.4
$ws.Range("I89").Value = 1481
$ws.Range("J89").Value = 1679.2858
$ws.Range("K89").Value = 7405
$ws.Range("L89").Value = 8396.429
$ws.Range("M89").Value = -1789
$ws.Range("N89").Value = -19628.429
$ws.Range("H97").Value = 7422
$ws.Range("J97").Value = 8777.5
$ws.Range("L97").Value = 26332.5
$ws.Range("N97").Value = -27324.5
$ws.Range("H141").Value = 4983.6924
$ws.Range("I141").Value = 5041.7754
$ws.Range("J141").Value = 4035
$ws.Range("K141").Value = 15125.3262
$ws.Range("L141").Value = 12105
$ws.Range("M141").Value = -9945.326200000001
$ws.Range("N141").Value = -22465

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 752.7436
$ws.Range("I61").Value = 503.1
$ws.Range("J61").Value = 1584.8889
$ws.Range("K61").Value = 503.1
$ws.Range("L61").Value = 1584.8889
$ws.Range("M61").Value = -291.1
$ws.Range("N61").Value = -2008.8889
$ws.Range("H122").Value = 2909.8333
$ws.Range("I122").Value = 1363
$ws.Range("J122").Value = 6003.5
$ws.Range("K122").Value = 4089
$ws.Range("L122").Value = 18010.5
$ws.Range("M122").Value = -1639
$ws.Range("N122").Value = -22910.5
$ws.Range("H132").Value = 1833.0769
$ws.Range("I132").Value = 1304
$ws.Range("J132").Value = 3326.9412
$ws.Range("K132").Value = 3912
$ws.Range("L132").Value = 9980.8236
$ws.Range("M132").Value = -1382
$ws.Range("N132").Value = -15040.8236
$ws.Range("H136").Value = 752.7436
$ws.Range("I136").Value = 503.1
$ws.Range("J136").Value = 1584.8889
$ws.Range("K136").Value = 1509.3
$ws.Range("L136").Value = 4754.6667
$ws.Range("M136").Value = 1040.7
$ws.Range("N136").Value = -9854.6667

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14234.75
$ws.Range("I20").Value = 2327
$ws.Range("J20").Value = 20188.625
$ws.Range("K20").Value = 2327
$ws.Range("L20").Value = 20188.625
$ws.Range("M20").Value = -2080
$ws.Range("N20").Value = -20682.625
$ws.Range("H80").Value = 136.3
$ws.Range("I80").Value = 77.416664
$ws.Range("J80").Value = 175.55556
$ws.Range("K80").Value = 77.416664
$ws.Range("L80").Value = 175.55556
$ws.Range("M80").Value = 920.583336
$ws.Range("N80").Value = -2171.55556
$ws.Range("H83").Value = 136.3
$ws.Range("I83").Value = 77.416664
$ws.Range("J83").Value = 175.55556
$ws.Range("K83").Value = 387.08332
$ws.Range("L83").Value = 877.7778000000001
$ws.Range("M83").Value = 4604.91668
$ws.Range("N83").Value = -10861.7778
$ws.Range("H94").Value = 447
$ws.Range("I94").Value = 546.36365
$ws.Range("J94").Value = 204.11111
$ws.Range("K94").Value = 546.36365
$ws.Range("L94").Value = 204.11111
$ws.Range("M94").Value = -95.36365000000001
$ws.Range("N94").Value = -1106.11111
$ws.Range("H105").Value = 1548.375
$ws.Range("I105").Value = 1404.7805
$ws.Range("J105").Value = 1940.8667
$ws.Range("K105").Value = 1404.7805
$ws.Range("L105").Value = 1940.8667
$ws.Range("M105").Value = 342.2194999999999
$ws.Range("N105").Value = -5434.8667
$ws.Range("H132").Value = 55555.555
$ws.Range("J132").Value = 55555.555
$ws.Range("L132").Value = 55555.555
$ws.Range("N132").Value = -65675.55499999999
$ws.Range("H134").Value = 1565.8073
$ws.Range("I134").Value = 922.3279
$ws.Range("J134").Value = 3350
$ws.Range("K134").Value = 2766.9837
$ws.Range("L134").Value = 10050
$ws.Range("M134").Value = -231.9836999999998
$ws.Range("N134").Value = -15120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 554
$ws.Range("I22").Value = 342.77274
$ws.Range("K22").Value = 342.77274
$ws.Range("M22").Value = 7.227260000000001
$ws.Range("H31").Value = 8066951.5
$ws.Range("I31").Value = 1271.4062
$ws.Range("J31").Value = 16670344
$ws.Range("K31").Value = 1271.4062
$ws.Range("L31").Value = 16670344
$ws.Range("M31").Value = -976.4061999999999
$ws.Range("N31").Value = -16670934
$ws.Range("H34").Value = 8066951.5
$ws.Range("I34").Value = 1271.4062
$ws.Range("J34").Value = 16670344
$ws.Range("K34").Value = 1271.4062
$ws.Range("L34").Value = 16670344
$ws.Range("M34").Value = -1069.4062
$ws.Range("N34").Value = -16670748
$ws.Range("H99").Value = 9096002
$ws.Range("I99").Value = 13336669
$ws.Range("K99").Value = 13336669
$ws.Range("M99").Value = -13335171
$ws.Range("H122").Value = 3376.6667
$ws.Range("I122").Value = 1696.6666
$ws.Range("J122").Value = 4216.6665
$ws.Range("K122").Value = 5089.9998
$ws.Range("L122").Value = 12649.9995
$ws.Range("M122").Value = -2639.9998
$ws.Range("N122").Value = -17549.9995
$ws.Range("H126").Value = 9096002
$ws.Range("I126").Value = 13336669
$ws.Range("K126").Value = 40010007
$ws.Range("M126").Value = -40007537
$ws.Range("H132").Value = 2971.0312
$ws.Range("I132").Value = 2459.2856
$ws.Range("K132").Value = 7377.8568
$ws.Range("M132").Value = -4847.8568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1384.1111
$ws.Range("I5").Value = 373.6842
$ws.Range("J5").Value = 2513.4119
$ws.Range("K5").Value = 1121.0526
$ws.Range("L5").Value = 7540.2357
$ws.Range("M5").Value = -1009.0526
$ws.Range("N5").Value = -7764.2357
$ws.Range("H107").Value = 21206.875
$ws.Range("I107").Value = 426.83334
$ws.Range("J107").Value = 33674.9
$ws.Range("K107").Value = 1280.50002
$ws.Range("L107").Value = 101024.7
$ws.Range("M107").Value = 639.4999800000001
$ws.Range("N107").Value = -104864.7
$ws.Range("H113").Value = 644.9286
$ws.Range("I113").Value = 650.5
$ws.Range("J113").Value = 631
$ws.Range("K113").Value = 1951.5
$ws.Range("L113").Value = 1893
$ws.Range("M113").Value = 218.5
$ws.Range("N113").Value = -6233
$ws.Range("H132").Value = 2055.5334
$ws.Range("I132").Value = 883.25
$ws.Range("J132").Value = 2481.818
$ws.Range("K132").Value = 7949.25
$ws.Range("L132").Value = 22336.362
$ws.Range("M132").Value = -5419.25
$ws.Range("N132").Value = -27396.362
$ws.Range("H135").Value = 1384.1111
$ws.Range("I135").Value = 373.6842
$ws.Range("J135").Value = 2513.4119
$ws.Range("K135").Value = 3363.1578
$ws.Range("L135").Value = 22620.7071
$ws.Range("M135").Value = -828.1578
$ws.Range("N135").Value = -27690.7071

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1502.2858
$ws.Range("I102").Value = 1060.0952
$ws.Range("J102").Value = 2165.5715
$ws.Range("K102").Value = 1060.0952
$ws.Range("L102").Value = 2165.5715
$ws.Range("M102").Value = 561.9048
$ws.Range("N102").Value = -5409.5715
$ws.Range("H126").Value = 3252.47
$ws.Range("I126").Value = 2979.9158
$ws.Range("J126").Value = 4583.1763
$ws.Range("K126").Value = 8939.7474
$ws.Range("L126").Value = 13749.5289
$ws.Range("M126").Value = -6469.7474
$ws.Range("N126").Value = -18689.5289
$ws.Range("H132").Value = 2047.9608
$ws.Range("I132").Value = 1313.7894
$ws.Range("J132").Value = 4194
$ws.Range("K132").Value = 3941.3682
$ws.Range("L132").Value = 12582
$ws.Range("M132").Value = -1411.3682
$ws.Range("N132").Value = -17642

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 41673470
$ws.Range("I40").Value = 111114480
$ws.Range("J40").Value = 8865.866
$ws.Range("K40").Value = 111114480
$ws.Range("L40").Value = 8865.866
$ws.Range("M40").Value = -111114344
$ws.Range("N40").Value = -9137.866
$ws.Range("H46").Value = 1372
$ws.Range("I46").Value = 1731.1428
$ws.Range("J46").Value = 1192.4286
$ws.Range("K46").Value = 1731.1428
$ws.Range("L46").Value = 1192.4286
$ws.Range("M46").Value = -1543.1428
$ws.Range("N46").Value = -1568.4286
$ws.Range("H132").Value = 9411.755999999999
$ws.Range("I132").Value = 11526.208
$ws.Range("J132").Value = 6995.2383
$ws.Range("K132").Value = 34578.624
$ws.Range("L132").Value = 20985.7149
$ws.Range("M132").Value = -32048.624
$ws.Range("N132").Value = -26045.7149
$ws.Range("H136").Value = 2045.6666
$ws.Range("I136").Value = 1143.4634
$ws.Range("J136").Value = 7330
$ws.Range("K136").Value = 3430.3902
$ws.Range("L136").Value = 21990
$ws.Range("M136").Value = -880.3902000000003
$ws.Range("N136").Value = -27090

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13335150
$ws.Range("I132").Value = 762.8421
$ws.Range("J132").Value = 55560708
$ws.Range("K132").Value = 2288.5263
$ws.Range("L132").Value = 166682124
$ws.Range("M132").Value = 241.4737
$ws.Range("N132").Value = -166687184
$ws.Range("H136").Value = 1868.5272
$ws.Range("I136").Value = 600.9211
$ws.Range("J136").Value = 4702
$ws.Range("K136").Value = 1802.7633
$ws.Range("L136").Value = 14106
$ws.Range("M136").Value = 747.2366999999999
$ws.Range("N136").Value = -19206
